$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force all touched cells to Text format first so numeric-looking
# strings (e.g. "0.530", "69.651.45") are preserved verbatim as text
# rather than being reinterpreted/reformatted as numbers by Excel.
$touchedCells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'E18', 'D19', 'E19', 'B20', 'C20', 'D20', 'E20', 'B21', 'C21', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D29', 'E29', 'D30', 'E30', 'D31', 'E31', 'E32', 'D33', 'E33', 'D34', 'E34', 'E35', 'D36', 'E36', 'D37', 'E37', 'D38', 'E38', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'B43', 'C43', 'D43', 'E43', 'B44', 'C44', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'B48', 'C48', 'D48', 'E48', 'B49', 'C49', 'D49', 'E49', 'E50', 'D51', 'E51')
foreach ($cellRef in $touchedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.651.45'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '3.726.29'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '612.79'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = '178.84'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('D7').Value = '3.722.40'
$ws.Range('E7').Value = '  -0.47%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.530'
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  -2.40%  '
$ws.Range('D11').Value = '6.57'
$ws.Range('E11').Value = '  +2.01%  '
$ws.Range('D12').Value = '0.480'
$ws.Range('E12').Value = '  -4.63%  '
$ws.Range('D13').Value = '39.79'
$ws.Range('E13').Value = '  -3.37%  '
$ws.Range('D14').Value = '0.0000254'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').Value = '4.348.17'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '3.726.55'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').Value = '69.741.34'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('E18').Value = '  -2.51%  '
$ws.Range('D19').Value = '7.49'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '501.91'
$ws.Range('E20').Value = '  -3.16%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = '16.32'
$ws.Range('E21').Value = '  -2.61%  '
$ws.Range('D22').Value = '9.12'
$ws.Range('E22').Value = '  -2.51%  '
$ws.Range('D23').Value = '0.719'
$ws.Range('E23').Value = '  -2.34%  '
$ws.Range('D24').Value = '2.61'
$ws.Range('E24').Value = '  +4.82%  '
$ws.Range('D25').Value = '86.18'
$ws.Range('E25').Value = '  -3.01%  '
$ws.Range('D26').Value = '11.30'
$ws.Range('E26').Value = '  +2.93%  '
$ws.Range('D27').Value = '12.97'
$ws.Range('E27').Value = '  -4.77%  '
$ws.Range('D28').Value = '0.0000136'
$ws.Range('E28').Value = '  +7.34%  '
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('D30').Value = '2.47'
$ws.Range('E30').Value = '  -1.80%  '
$ws.Range('D31').Value = '2.91'
$ws.Range('E31').Value = '  +1.66%  '
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').Value = '30.48'
$ws.Range('E33').Value = '  -3.57%  '
$ws.Range('D34').Value = '0.113'
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('D36').Value = '1.04'
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('D37').Value = '6.11'
$ws.Range('E37').Value = '  -2.31%  '
$ws.Range('D38').Value = '0.349'
$ws.Range('E38').Value = '  +2.19%  '
$ws.Range('E39').Value = '  +5.87%  '
$ws.Range('D40').Value = '3.09'
$ws.Range('E40').Value = '  +12.72%  '
$ws.Range('D41').Value = '2.06'
$ws.Range('E41').Value = '  -5.82%  '
$ws.Range('D42').Value = '45.84'
$ws.Range('E42').Value = '  +3.14%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '49.65'
$ws.Range('E43').Value = '  -3.51%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').Value = '436.58'
$ws.Range('E44').Value = '  +3.03%  '
$ws.Range('D45').Value = '8.56'
$ws.Range('E45').Value = '  -3.57%  '
$ws.Range('D46').Value = '2.951.52'
$ws.Range('E46').Value = '  -3.96%  '
$ws.Range('D47').Value = '0.0362'
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '139.15'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '27.15'
$ws.Range('E49').Value = '  -2.83%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').Value = '2.46'
$ws.Range('E51').Value = '  -2.39%  '
